$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 35968.234
$ws.Cells.Item(64, 9).Value = 113466.664
$ws.Cells.Item(64, 10).Value = 2754.6191
$ws.Cells.Item(64, 11).Value = 113466.664
$ws.Cells.Item(64, 12).Value = 2754.6191
$ws.Cells.Item(64, 13).Value = -113218.664
$ws.Cells.Item(64, 14).Value = -3250.6191
# Row 67
$ws.Cells.Item(67, 8).Value = 35968.234
$ws.Cells.Item(67, 9).Value = 113466.664
$ws.Cells.Item(67, 10).Value = 2754.6191
$ws.Cells.Item(67, 11).Value = 113466.664
$ws.Cells.Item(67, 12).Value = 2754.6191
$ws.Cells.Item(67, 13).Value = -112608.664
$ws.Cells.Item(67, 14).Value = -4470.6191
# Row 74
$ws.Cells.Item(74, 8).Value = 3638.3845
$ws.Cells.Item(74, 9).Value = 4033.3333
$ws.Cells.Item(74, 10).Value = 3299.8572
$ws.Cells.Item(74, 11).Value = 4033.3333
$ws.Cells.Item(74, 12).Value = 3299.8572
$ws.Cells.Item(74, 13).Value = -3097.3333
$ws.Cells.Item(74, 14).Value = -5171.8572
# Row 76
$ws.Cells.Item(76, 8).Value = 4825.75
$ws.Cells.Item(76, 9).Value = 5334.3335
$ws.Cells.Item(76, 10).Value = 3300
$ws.Cells.Item(76, 11).Value = 5334.3335
$ws.Cells.Item(76, 12).Value = 3300
$ws.Cells.Item(76, 13).Value = -5019.3335
$ws.Cells.Item(76, 14).Value = -3930
# Row 77
$ws.Cells.Item(77, 8).Value = 3638.3845
$ws.Cells.Item(77, 9).Value = 4033.3333
$ws.Cells.Item(77, 10).Value = 3299.8572
$ws.Cells.Item(77, 11).Value = 20166.6665
$ws.Cells.Item(77, 12).Value = 16499.286
$ws.Cells.Item(77, 13).Value = -15486.6665
$ws.Cells.Item(77, 14).Value = -25859.286
# Row 79
$ws.Cells.Item(79, 8).Value = 4825.75
$ws.Cells.Item(79, 9).Value = 5334.3335
$ws.Cells.Item(79, 10).Value = 3300
$ws.Cells.Item(79, 11).Value = 5334.3335
$ws.Cells.Item(79, 12).Value = 3300
$ws.Cells.Item(79, 13).Value = -4242.3335
$ws.Cells.Item(79, 14).Value = -5484
# Row 104
$ws.Cells.Item(104, 8).Value = 3026.7778
$ws.Cells.Item(104, 9).Value = 191.5
$ws.Cells.Item(104, 11).Value = 574.5
$ws.Cells.Item(104, 13).Value = 1172.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 135
$ws.Cells.Item(135, 8).Value = 35498
$ws.Cells.Item(135, 10).Value = 35498
$ws.Cells.Item(135, 12).Value = 35498
$ws.Cells.Item(135, 14).Value = -45638
# Row 138
$ws.Cells.Item(138, 8).Value = 17049.285
$ws.Cells.Item(138, 10).Value = 17049.285
$ws.Cells.Item(138, 12).Value = 17049.285
$ws.Cells.Item(138, 14).Value = -27329.285
# Row 140
$ws.Cells.Item(140, 8).Value = 31717.4
$ws.Cells.Item(140, 10).Value = 31717.4
$ws.Cells.Item(140, 12).Value = 31717.4
$ws.Cells.Item(140, 14).Value = -42077.4

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Cells.Item(6, 8).Value = 2944144.5
$ws.Cells.Item(6, 9).Value = 20000000
$ws.Cells.Item(6, 11).Value = 20000000
$ws.Cells.Item(6, 13).Value = -19999887
# Row 7
$ws.Cells.Item(7, 8).Value = 392.7857
$ws.Cells.Item(7, 9).Value = 293.85715
$ws.Cells.Item(7, 10).Value = 491.7143
$ws.Cells.Item(7, 11).Value = 293.85715
$ws.Cells.Item(7, 12).Value = 491.7143
$ws.Cells.Item(7, 13).Value = -180.85715
$ws.Cells.Item(7, 14).Value = -717.7143
# Row 17
$ws.Cells.Item(17, 8).Value = 19500
$ws.Cells.Item(17, 9).Value = 19000
$ws.Cells.Item(17, 10).Value = 20000
$ws.Cells.Item(17, 11).Value = 19000
$ws.Cells.Item(17, 12).Value = 20000
$ws.Cells.Item(17, 13).Value = -18826
$ws.Cells.Item(17, 14).Value = -20348
# Row 31
$ws.Cells.Item(31, 8).Value = 5436.552
$ws.Cells.Item(31, 9).Value = 2281.8096
$ws.Cells.Item(31, 10).Value = 7227.081
$ws.Cells.Item(31, 11).Value = 2281.8096
$ws.Cells.Item(31, 12).Value = 7227.081
$ws.Cells.Item(31, 13).Value = -1986.8096
$ws.Cells.Item(31, 14).Value = -7817.081
# Row 34
$ws.Cells.Item(34, 8).Value = 5436.552
$ws.Cells.Item(34, 9).Value = 2281.8096
$ws.Cells.Item(34, 10).Value = 7227.081
$ws.Cells.Item(34, 11).Value = 2281.8096
$ws.Cells.Item(34, 12).Value = 7227.081
$ws.Cells.Item(34, 13).Value = -2079.8096
$ws.Cells.Item(34, 14).Value = -7631.081
# Row 41
$ws.Cells.Item(41, 8).Value = 25456.445
$ws.Cells.Item(41, 9).Value = 6666.6665
$ws.Cells.Item(41, 10).Value = 34851.332
$ws.Cells.Item(41, 11).Value = 6666.6665
$ws.Cells.Item(41, 12).Value = 34851.332
$ws.Cells.Item(41, 13).Value = -6238.6665
$ws.Cells.Item(41, 14).Value = -35707.332
# Row 50
$ws.Cells.Item(50, 8).Value = 37959.75
$ws.Cells.Item(50, 10).Value = 37959.75
$ws.Cells.Item(50, 12).Value = 37959.75
$ws.Cells.Item(50, 14).Value = -39209.75
# Row 59
$ws.Cells.Item(59, 8).Value = 38858
$ws.Cells.Item(59, 10).Value = 40322.5
$ws.Cells.Item(59, 12).Value = 40322.5
$ws.Cells.Item(59, 14).Value = -42612.5
# Row 60
$ws.Cells.Item(60, 8).Value = 23603
$ws.Cells.Item(60, 10).Value = 23603
$ws.Cells.Item(60, 12).Value = 23603
$ws.Cells.Item(60, 14).Value = -24625
# Row 62
$ws.Cells.Item(62, 8).Value = 2816.6667
$ws.Cells.Item(62, 9).Value = 2875
$ws.Cells.Item(62, 10).Value = 2787.5
$ws.Cells.Item(62, 11).Value = 2875
$ws.Cells.Item(62, 12).Value = 2787.5
$ws.Cells.Item(62, 13).Value = -2251
$ws.Cells.Item(62, 14).Value = -4035.5
# Row 65
$ws.Cells.Item(65, 8).Value = 2816.6667
$ws.Cells.Item(65, 9).Value = 2875
$ws.Cells.Item(65, 10).Value = 2787.5
$ws.Cells.Item(65, 11).Value = 14375
$ws.Cells.Item(65, 12).Value = 13937.5
$ws.Cells.Item(65, 13).Value = -11255
$ws.Cells.Item(65, 14).Value = -20177.5
# Row 74
$ws.Cells.Item(74, 8).Value = 13000
$ws.Cells.Item(74, 10).Value = 13000
$ws.Cells.Item(74, 12).Value = 13000
$ws.Cells.Item(74, 14).Value = -14748
# Row 77
$ws.Cells.Item(77, 8).Value = 13000
$ws.Cells.Item(77, 10).Value = 13000
$ws.Cells.Item(77, 12).Value = 39000
$ws.Cells.Item(77, 14).Value = -47736
# Row 80
$ws.Cells.Item(80, 8).Value = 30459
$ws.Cells.Item(80, 10).Value = 30459
$ws.Cells.Item(80, 12).Value = 30459
$ws.Cells.Item(80, 14).Value = -32705
# Row 83
$ws.Cells.Item(83, 8).Value = 30459
$ws.Cells.Item(83, 10).Value = 30459
$ws.Cells.Item(83, 12).Value = 91377
$ws.Cells.Item(83, 14).Value = -102609

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).Value = ""
# Row 70
$ws.Cells.Item(70, 8).Value = 5238.8716
$ws.Cells.Item(70, 9).Value = 5103.8667
$ws.Cells.Item(70, 10).Value = 5688.8887
$ws.Cells.Item(70, 11).Value = 5103.8667
$ws.Cells.Item(70, 12).Value = 5688.8887
$ws.Cells.Item(70, 13).Value = -4833.8667
$ws.Cells.Item(70, 14).Value = -6228.8887
# Row 73
$ws.Cells.Item(73, 8).Value = 5238.8716
$ws.Cells.Item(73, 9).Value = 5103.8667
$ws.Cells.Item(73, 10).Value = 5688.8887
$ws.Cells.Item(73, 11).Value = 5103.8667
$ws.Cells.Item(73, 12).Value = 5688.8887
$ws.Cells.Item(73, 13).Value = -4167.8667
$ws.Cells.Item(73, 14).Value = -7560.8887
# Row 80
$ws.Cells.Item(80, 8).Value = 212045.83
$ws.Cells.Item(80, 10).Value = 3069.2307
$ws.Cells.Item(80, 12).Value = 3069.2307
$ws.Cells.Item(80, 14).Value = -5065.2307
# Row 81
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).Value = ""
# Row 83
$ws.Cells.Item(83, 8).Value = 212045.83
$ws.Cells.Item(83, 10).Value = 3069.2307
$ws.Cells.Item(83, 12).Value = 15346.1535
$ws.Cells.Item(83, 14).Value = -25330.1535
# Row 84
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).Value = ""
# Row 126
$ws.Cells.Item(126, 8).Value = 5778.593
$ws.Cells.Item(126, 9).Value = 9051.571
$ws.Cells.Item(126, 10).Value = 2253.8462
$ws.Cells.Item(126, 11).Value = 27154.713
$ws.Cells.Item(126, 12).Value = 6761.5386
$ws.Cells.Item(126, 13).Value = -24684.713
$ws.Cells.Item(126, 14).Value = -11701.5386

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 9
$ws.Cells.Item(9, 8).Value = 10700
$ws.Cells.Item(9, 9).Value = 400
$ws.Cells.Item(9, 11).Value = 400
$ws.Cells.Item(9, 13).Value = -176
